# Built reusable QuestionForm component and remove redundant code from add and edit question
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: set Correct answers (D) to 0 and Incorrect answers (E) to 2
# for rows where a retake/re-evaluation happened.
$rowsToUpdate = @(8, 11, 13, 14, 15, 16, 21, 27, 30, 31)
foreach ($r in $rowsToUpdate) {
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 2
}

# Append new quiz attempt rows (32-35)
$ws.Cells.Item(32, 1).Value = "General Knowledge"
$ws.Cells.Item(32, 2).Value = "Sanjib Roy"
$ws.Cells.Item(32, 3).Value = "sanjibroysnjsnsj0098@gmail.com"
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 2

$ws.Cells.Item(33, 1).Value = "General Knowledge"
$ws.Cells.Item(33, 2).Value = "Sanjib Roy"
$ws.Cells.Item(33, 3).Value = "sanjibroy53110098@gmail.com"
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 2

$ws.Cells.Item(34, 1).Value = "General Knowledge"
$ws.Cells.Item(34, 2).Value = "Sanjib Roy"
$ws.Cells.Item(34, 3).Value = "sanjibroy00444698@gmail.com"
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 2

$ws.Cells.Item(35, 1).Value = "General Knowledge"
$ws.Cells.Item(35, 2).Value = "Sanjib Roy"
$ws.Cells.Item(35, 3).Value = "san4646546jibroy0098@gmail.com"
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 1
